# Update "Datos actualizados" timestamp and refresh COVID-19 country stats
# (commit: "Update countries & provincias Spain")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 10:22"

# Row 4 (Estados Unidos)
$ws.Cells.Item(4, 2).Value = 6389413
$ws.Cells.Item(4, 3).Value = 356
$ws.Cells.Item(4, 4).Value = 3636272
$ws.Cells.Item(4, 5).Value = 2561009
$ws.Cells.Item(4, 7).Value = 21
$ws.Cells.Item(4, 8).Value = 192132

# Row 6 (Peru)
$ws.Cells.Item(6, 2).Value = 4027718
$ws.Cells.Item(6, 3).Value = 7479
$ws.Cells.Item(6, 4).Value = 3107453
$ws.Cells.Item(6, 5).Value = 850597
$ws.Cells.Item(6, 7).Value = 33
$ws.Cells.Item(6, 8).Value = 69668

# Row 7 (Colombia)
$ws.Cells.Item(7, 2).Value = 1020310
$ws.Cells.Item(7, 3).Value = 5205
$ws.Cells.Item(7, 4).Value = 838126
$ws.Cells.Item(7, 5).Value = 164425
$ws.Cells.Item(7, 7).Value = 110
$ws.Cells.Item(7, 8).Value = 17759

# Row 25
$ws.Cells.Item(25, 2).Value = 234570
$ws.Cells.Item(25, 3).Value = 2529
$ws.Cells.Item(25, 4).Value = 161668
$ws.Cells.Item(25, 5).Value = 69112
$ws.Cells.Item(25, 7).Value = 53
$ws.Cells.Item(25, 8).Value = 3790

# Row 48
$ws.Cells.Item(48, 4).Value = 52346
$ws.Cells.Item(48, 5).Value = 15374

# Row 52
$ws.Cells.Item(52, 2).Value = 56982
$ws.Cells.Item(52, 3).Value = 34
$ws.Cells.Item(52, 5).Value = 781

# Row 66
$ws.Cells.Item(66, 2).Value = 38324
$ws.Cells.Item(66, 3).Value = 20
$ws.Cells.Item(66, 4).Value = 30082
$ws.Cells.Item(66, 5).Value = 6833

# Row 74
$ws.Cells.Item(74, 2).Value = 26206
$ws.Cells.Item(74, 3).Value = 107
$ws.Cells.Item(74, 4).Value = 15587
$ws.Cells.Item(74, 5).Value = 9867

# Row 137
$ws.Cells.Item(137, 2).Value = 2491
$ws.Cells.Item(137, 3).Value = 35
$ws.Cells.Item(137, 4).Value = 2165
$ws.Cells.Item(137, 5).Value = 262

# Row 155
$ws.Cells.Item(155, 2).Value = 1425
$ws.Cells.Item(155, 3).Value = 9
$ws.Cells.Item(155, 5).Value = 203

# Row 165
$ws.Cells.Item(165, 4).Value = 786
$ws.Cells.Item(165, 5).Value = 228
